# Updated cryptos list (price / 1h volume refresh + a few rank swaps).
# Values are set with a leading apostrophe to force text (so numeric-looking
# strings like "609.01" or "1.00" are stored verbatim instead of being
# parsed into numbers / losing trailing zeros), then the style is reset to
# "Normal" so the quote-prefix formatting doesn't linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.044.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.73%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.717.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.75%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'609.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.37%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'189.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.76%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.37%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.718"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.84%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'58.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +7.07%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0000289"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.42%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'10.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.56%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.303.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.56%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.715.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.84%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.34%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'19.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.48%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -1.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'68.813.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.69%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'410.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.69%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.29%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'89.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.14%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.66%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'12.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.16%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.95%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'6.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.23%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.16%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.61%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'33.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.75%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -8.79%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'12.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.15%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'InjectiveProtocol"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'46.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.57%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'Bittensor"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'641.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.85%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'Hedera"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.67%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'65.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.95%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.411"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0₃0821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -12.36%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.14%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.07%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.140"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.32%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -3.02%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.42%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.80%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +2.18%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.854.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.46%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'THORChain"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'9.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.84%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'WEMIXToken"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'2.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.31%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.29%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'141.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.85%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -22.94%  "
$ws.Range("E51").Style = "Normal"
